# Reduce the VRP instance from 4 intermediate nodes (N1-N4) down to 3
# (N1-N3): drop the N4 row/column on every sheet and refresh the
# remaining numbers + time windows for the smaller verification instance.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "demand": drop the N4 row, then rewrite the surviving rows
# ---------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("demand")
$wsDemand.Rows.Item(6).Delete()

# Row 2 - DepotStart
$wsDemand.Range("F2").Value = 3600

# Row 3 - N1
$wsDemand.Range("B3").Value = 1
$wsDemand.Range("C3").Value = 361
$wsDemand.Range("D3").Value = "14:00 - 14:30"
$wsDemand.Range("E3").Value = 21600
$wsDemand.Range("F3").Value = 23400

# Row 4 - N2
$wsDemand.Range("B4").Value = 1
$wsDemand.Range("C4").Value = 440
$wsDemand.Range("D4").Value = "13:00 - 13:30"
$wsDemand.Range("E4").Value = 18000
$wsDemand.Range("F4").Value = 19800

# Row 5 - N3
$wsDemand.Range("B5").Value = 5
$wsDemand.Range("C5").Value = 429
$wsDemand.Range("D5").Value = "9:00 - 9:30"
$wsDemand.Range("E5").Value = 3600
$wsDemand.Range("F5").Value = 5400

# Row 6 - DepotEnd (was row 7 before the delete pulled it up)
$wsDemand.Range("B6").Value = 0
$wsDemand.Range("C6").Value = 0
$wsDemand.Range("E6").Value = 1
$wsDemand.Range("F6").Value = 46800

# ---------------------------------------------------------------
# Sheet "Distance": drop the N4 row + column, then rewrite values
# ---------------------------------------------------------------
$wsDist = $wb.Worksheets.Item("Distance")
$wsDist.Rows.Item(6).Delete()
$wsDist.Columns.Item(7).Delete()
$wsDist.Range("F1").Value = "DepotEnd"

$distValues = @(
    @(0, 6, 13, 11, 0),
    @(6, 0, 19, 5, 6),
    @(13, 19, 0, 14, 13),
    @(11, 5, 14, 0, 11),
    @(0, 6, 13, 11, 0)
)
for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $wsDist.Cells.Item($r + 2, $c + 2).Value = $distValues[$r][$c]
    }
}

# ---------------------------------------------------------------
# Sheet "TravelTime": drop the N4 row + column, then rewrite values
# ---------------------------------------------------------------
$wsTime = $wb.Worksheets.Item("TravelTime")
$wsTime.Rows.Item(6).Delete()
$wsTime.Columns.Item(7).Delete()
$wsTime.Range("F1").Value = "DepotEnd"

$timeValues = @(
    @(0, 720, 1560, 1320, 0),
    @(0, 0, 2280, 600, 720),
    @(0, 2280, 0, 1680, 1560),
    @(0, 600, 1680, 0, 1320),
    @(0, 0, 0, 0, 0)
)
for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $wsTime.Cells.Item($r + 2, $c + 2).Value = $timeValues[$r][$c]
    }
}
